$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the middle "1912" data row (old row 17). Excel shifts row 18
# (and everything below it, incl. the signature rows) up by one, which
# also auto-updates dimension / mergeCells / shared strings.
$ws.Rows(17).Delete()

# Update the remaining figures to their new values.
$ws.Range("E11").Value2 = 240000       # VALOR MORA (total)
$ws.Range("F13").Value2 = 2            # Cant. Periodos

$ws.Range("E16").Value2 = "1911"       # Periodo Mora for the first data row
$ws.Range("F16").Value2 = 120000       # Valor Mora for period 1911 row
$ws.Range("G16").Value2 = 3120000      # Salario Basico for period 1911 row

$ws.Range("E17").Value2 = "1912"       # Periodo Mora for the remaining data row (now row 17)
$ws.Range("G17").Value2 = 3120000      # Salario Basico for period 1912 row (now row 17)
